$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Range("S2").Value = "Pass"
$ws.Range("S3").Value = "Pass"
$ws.Range("S6").Value = "Pass"
$ws.Range("S2:S10").Select()
$ws.Cells.Item(10, 19).Activate()
